$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edits -------------------------------------------------------
# A2: joao126 -> joao1211
$ws.Range("A2").Value = "joao1211"

# Row 3 (Pedro999 test scenario) cleared out - only B3 keeps its style
$ws.Range("A3:L3").ClearContents()

# --- Selection ----------------------------------------------------------
$ws.Range("B10").Select()
